$d = $word.ActiveDocument

# 1. Insert a new standalone paragraph "What is a frontend, what is a backend?"
#    right before the standalone "HTML" paragraph (currently paragraph 3).
$pHtml = $d.Paragraphs(3)
$pHtml.Range.InsertParagraphBefore()
$d.Paragraphs(3).Range.Text = "What is a frontend, what is a backend?"

# 2. Remove the two list items "Lists & Tables" and
#    "Images (Talk about </> and <> </> tags)" that used to follow "Headings".
#    (Indices shifted by +1 because of the insertion above: 16 -> 17, 17 -> 18.)
$d.Paragraphs(18).Range.Delete()
$d.Paragraphs(17).Range.Delete()

# 3. Replace the two trailing empty paragraphs with a single paragraph of text.
$count = $d.Paragraphs.Count
$d.Paragraphs($count - 1).Range.Delete()
$count = $d.Paragraphs.Count
$d.Paragraphs($count).Range.Text = "More resources (w3schools, Mozilla docs, online courses)"

# 4. Merge the "Node" / " & " runs into a single run (cosmetic run-consolidation,
#    matching what Word does when that text is touched during editing).
$d.Content.Find.Execute("Node & ", $false, $false, $false, $false, $false, $true,
                         1, $false, "Node & ", 2) | Out-Null
